# Add a "slug" column to the "news" sheet (between "title" and "excerpt"),
# populate it for the existing news row, and leave the "news" sheet active
# with cell C3 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("news")
$ws.Activate()

# Shift the existing columns C:G one column to the right (into D:H) to make
# room for the new "slug" column at C. Work from right to left so a value
# isn't overwritten before it has been copied onward.
$ws.Range("H1").Value2 = $ws.Range("G1").Value2
$ws.Range("H2").Value2 = $ws.Range("G2").Value2

$ws.Range("G1").Value2 = $ws.Range("F1").Value2
$ws.Range("G2").Value2 = $ws.Range("F2").Value2

$ws.Range("F1").Value2 = $ws.Range("E1").Value2
$ws.Range("F2").Value2 = $ws.Range("E2").Value2

$ws.Range("E1").Value2 = $ws.Range("D1").Value2
$ws.Range("E2").Value2 = $ws.Range("D2").Value2

$ws.Range("D1").Value2 = $ws.Range("C1").Value2
$ws.Range("D2").Value2 = $ws.Range("C2").Value2

# Populate the new "slug" column.
$ws.Range("C1").Value2 = "slug"
$ws.Range("C2").Value2 = "ciampitti-lab-new-chapter"

# Match the saved selection/active sheet state.
$ws.Range("C3").Select()
